$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 00:10"

# Swap country names for rows 54/55: Bielorrusia now ranks above Honduras
$ws.Range("A54").Value = "Bielorrusia"
$ws.Range("A55").Value = "Honduras"

# Refresh per-country statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes)
$ws.Range("B4").Value = 7988371
$ws.Range("C4").Value = 40081
$ws.Range("D4").Value = 5124756
$ws.Range("E4").Value = 2643931
$ws.Range("G4").Value = 314
$ws.Range("H4").Value = 219684

$ws.Range("D6").Value = 4470163
$ws.Range("E6").Value = 474935
$ws.Range("G6").Value = 252
$ws.Range("H6").Value = 150488

$ws.Range("B34").Value = 147033
$ws.Range("C34").Value = 205
$ws.Range("D34").Value = 128134
$ws.Range("E34").Value = 6708
$ws.Range("G34").Value = 3
$ws.Range("H34").Value = 12191

$ws.Range("B53").Value = 84295
$ws.Range("C53").Value = 866
$ws.Range("D53").Value = 38316
$ws.Range("E53").Value = 44692
$ws.Range("G53").Value = 10
$ws.Range("H53").Value = 1287

$ws.Range("B54").Value = 83534
$ws.Range("C54").Value = 1063
$ws.Range("D54").Value = 77220
$ws.Range("E54").Value = 5418
$ws.Range("G54").Value = 11
$ws.Range("H54").Value = 896

$ws.Range("B55").Value = 83146
$ws.Range("C55").Value = 594
$ws.Range("D55").Value = 31931
$ws.Range("E55").Value = 48711
$ws.Range("G55").Value = 12
$ws.Range("H55").Value = 2504

$ws.Range("B86").Value = 24402
$ws.Range("C86").Value = 83
$ws.Range("D86").Value = 15847
$ws.Range("E86").Value = 7663
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 892

$ws.Range("B96").Value = 15458
$ws.Range("C96").Value = 43
$ws.Range("D96").Value = 14599
$ws.Range("E96").Value = 522

$ws.Range("B100").Value = 13673
$ws.Range("C100").Value = 3
$ws.Range("E100").Value = 6073

$ws.Range("B115").Value = 8011
$ws.Range("C115").Value = 1
$ws.Range("D115").Value = 6504
$ws.Range("E115").Value = 1277

$ws.Range("B146").Value = 3469
$ws.Range("C146").Value = 64
$ws.Range("D146").Value = 2318
$ws.Range("E146").Value = 1048
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 103

$ws.Range("B159").Value = 2052
$ws.Range("C159").Value = 1
$ws.Range("E159").Value = 128

$ws.Range("B161").Value = 1940
$ws.Range("C161").Value = 5
$ws.Range("D161").Value = 1457
$ws.Range("E161").Value = 434

